{"js": "// \"Fixat lite med bossen.\"\n//\n// 1. Merge the four title runs (\"Iteration 1\" / \" - \" / \"Elaboration (v.16\" / \")\")\n//    into a single run reading \"Iteration 1 - Elaboration (v.16)\".\n// 2. Delete the trailing \"Total tid\" row from the time-report table.\n// 3. Move the _GoBack bookmark from the very last paragraph (end of document)\n//    up to the now-empty paragraph that immediately follows the table.\n\n// --- 1. Collapse the title's four runs into one -----------------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nconst titleRange = titlePara.getRange();\ntitleRange.insertText(\"Iteration 1 \u2013 Elaboration (v.16)\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2. Delete the \"Total tid\" row (last row of the first table) -----------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.items[rows.items.length - 1].delete();\nawait context.sync();\n\n// --- 3. Relocate the _GoBack bookmark ---------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst afterTableRange = table.getRange(\"End\");\nconst afterTablePara = afterTableRange.paragraphs.getFirst();\nafterTablePara.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# \"Fixat lite med bossen.\"\n#\n# 1. Merge the four title runs (\"Iteration 1\" / \" - \" / \"Elaboration (v.16\" / \")\")\n#    into a single run reading \"Iteration 1 - Elaboration (v.16)\".\n# 2. Delete the trailing \"Total tid\" row from the time-report table.\n# 3. Move the _GoBack bookmark from the very last paragraph (end of document)\n#    up to the now-empty paragraph that immediately follows the table.\n\n$d = $word.ActiveDocument\n\n# --- 1. Collapse the title's four runs into one -----------------------------\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n$titleRange.MoveEnd(1, -1) | Out-Null                 # exclude the paragraph mark\n$firstRun = $titleRange.Duplicate\n$firstRun.End = $firstRun.Start + 11                  # \"Iteration 1\" (11 chars)\n$remainder = $titleRange.Duplicate\n$remainder.Start = $firstRun.End                      # \" - Elaboration (v.16)\"\n$remainderText = $remainder.Text\n\n$firstRun.Collapse(0)                                 # wdCollapseEnd\n$firstRun.InsertAfter($remainderText)\n$remainder.Delete()\n\n# --- 2. Delete the \"Total tid\" row (last row of the first table) -----------\n$table = $d.Tables.Item(1)\n$table.Rows.Last.Delete()\n\n# --- 3. Relocate the _GoBack bookmark ---------------------------------------\n$existing = $d.Bookmarks.Item(\"_GoBack\")\n$existing.Delete()\n\n$afterTableRange = $table.Range.Next(4, 1)            # wdParagraph = 4\n$afterTablePara = $afterTableRange.Paragraphs.Item(1)\n$d.Bookmarks.Add(\"_GoBack\", $afterTablePara.Range)\n"}
